$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 is the ferrite bead line item (FB1..FB7). Correct the part number
# (Value / Device columns) and its Description to match the new part.
$ws.Range("E9").Value = "HH-1T2012-601"
$ws.Range("H9").Value = "FERRITE CHIP 600 OHM 2500MA 0805"
